# This script re-orders the observation rows (2-10) on the active sheet,
# moving each row's full set of attributes (columns C..AY) to a new row
# position, and assigning new "Id" (column A) and "Taxonsorteringsordning"
# (column B) values to each resulting row, per the source diff.
#
# Because every row 2-10 is simultaneously a copy source and a copy
# target (a full permutation made of two cycles), we first stash a
# faithful copy of each source row in a scratch area far below the
# data, then paste each stashed row into its final destination. Using
# Copy + PasteSpecial(xlPasteValues) (rather than .Value/.Value2
# assignment) keeps every cell's original data type intact -- in
# particular it prevents Excel from "helpfully" reinterpreting the
# literal text dates/times in columns Y/Z/AA/AB (e.g. "2023-09-14",
# "12:01") as real date/time serial numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues

$sourceRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10)
$scratchOffset = 500

# Step 1: stash a values-only copy of each source row into scratch rows
# far below the real data (row + 500), so later writes to the real rows
# can't clobber data we still need to read.
foreach ($r in $sourceRows) {
    $ws.Range("A" + $r + ":AY" + $r).Copy()
    $scratchRow = $r + $scratchOffset
    $ws.Range("A" + $scratchRow).PasteSpecial($xlPasteValues)
}
$excel.CutCopyMode = $false

# Step 2: mapping of destination row -> the row whose stashed content
# should populate it, plus the new Id (A) / Taxonsorteringsordning (B)
# values for that destination row.
$rowMap = @{
    2  = 9
    3  = 4
    4  = 6
    5  = 8
    6  = 5
    7  = 2
    8  = 10
    9  = 7
    10 = 3
}
$newA = @{
    2  = 112086407
    3  = 112086079
    4  = 112086235
    5  = 112085339
    6  = 112085469
    7  = 112085285
    8  = 112085668
    9  = 112085312
    10 = 112086207
}
$newB = @{
    2  = 56446
    3  = 90792
    4  = 93539
    5  = 88953
    6  = 56446
    7  = 89539
    8  = 89539
    9  = 89539
    10 = 56446
}

# Step 3: copy each stashed row into its destination row, then overwrite
# the Id / Taxonsorteringsordning cells with their new values. The
# destination row is cleared first because PasteSpecial only overwrites
# cells that actually hold a value in the copied source range, and would
# otherwise leave behind stray leftover cells (e.g. a destination row
# that used to have a "Publik kommentar" / AC value but whose new
# content shouldn't have one).
foreach ($destRow in $sourceRows) {
    $ws.Range("A" + $destRow + ":AY" + $destRow).ClearContents()

    $srcRow = $rowMap[$destRow] + $scratchOffset
    $ws.Range("A" + $srcRow + ":AY" + $srcRow).Copy()
    $ws.Range("A" + $destRow).PasteSpecial($xlPasteValues)

    $ws.Range("A" + $destRow).Value2 = $newA[$destRow]
    $ws.Range("B" + $destRow).Value2 = $newB[$destRow]
}
$excel.CutCopyMode = $false

# Step 4: remove the scratch rows so the sheet's used range / dimension
# goes back to its original extent.
foreach ($r in $sourceRows) {
    $scratchRow = $r + $scratchOffset
    $ws.Range("A" + $scratchRow + ":AY" + $scratchRow).Clear()
}
